$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29/30 swap: Toncoin moved up in ranking, swapping places with Cosmos ---
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.12"
$ws.Range("E29").Value = "  +0.94%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.48"
$ws.Range("E30").Value = "  -0.18%  "

# --- Price / Volume(1h) updates ---
$ws.Range("D2").Value = "42.054.43"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.248.36"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'306.56"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'96.40"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.487"
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").Value = "'34.71"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "2.597.17"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "'14.45"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "2.258.27"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "'0.778"
$ws.Range("E17").Value = "  -2.37%  "
$ws.Range("D18").Value = "41.920.45"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "'12.14"
$ws.Range("E19").Value = "  -3.31%  "
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").Value = "'5.92"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "'67.12"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").Value = "'235.22"
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("D24").Value = "'2.57"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'23.32"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "'36.87"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D31").Value = "'165.06"
$ws.Range("E31").Value = "  +3.26%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("E35").Value = "  +3.47%  "
$ws.Range("D36").Value = "'0.0719"
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "'0.103"
$ws.Range("E39").Value = "  -3.04%  "
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").Value = "1.943.35"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").Value = "'2.19"
$ws.Range("E44").Value = "  -9.01%  "
$ws.Range("D45").Value = "'18.43"
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("E46").Value = "  -2.93%  "
$ws.Range("D47").Value = "'9.67"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").Value = "'53.58"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").Value = "2.469.01"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("D50").Value = "'71.23"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "'91.07"
$ws.Range("E51").Value = "  -0.87%  "
